$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Range("B38").Value = 6802967
$ws.Range("E38").Value = "Viktoria Plzen"
$ws.Range("F38").Value = "Sigma Olomouc"
$ws.Range("H38").Value = 1
$ws.Range("I38").Value = 2
$ws.Range("J38").Value = 1
$ws.Range("L38").Value = 1.55
$ws.Range("M38").Value = 3.8
$ws.Range("N38").Value = 5.25
$ws.Range("O38").Value = 1.45
$ws.Range("P38").Value = 4.2
$ws.Range("Q38").Value = 5.75
$ws.Range("R38").Value = -1
$ws.Range("S38").Value = 1.825
$ws.Range("T38").Value = 2.025
$ws.Range("U38").Value = 3
$ws.Range("V38").Value = 2
$ws.Range("W38").Value = 1.85
$ws.Range("X38").Value = 0.45
$ws.Range("AA38").Value = 0
$ws.Range("AB38").Value = 0
$ws.Range("AC38").Value = 0
$ws.Range("AD38").Value = 0

# Row 39
$ws.Range("B39").Value = 6802969
$ws.Range("E39").Value = "Hradec Kralove"
$ws.Range("F39").Value = "FC Trinity Zlin"
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 1.85
$ws.Range("M39").Value = 3.4
$ws.Range("N39").Value = 3.75
$ws.Range("O39").Value = 1.75
$ws.Range("P39").Value = 3.6
$ws.Range("Q39").Value = 4.2
$ws.Range("R39").Value = -0.75
$ws.Range("S39").Value = 2.025
$ws.Range("T39").Value = 1.825
$ws.Range("U39").Value = 2.75
$ws.Range("V39").Value = 2.025
$ws.Range("W39").Value = 1.825
$ws.Range("X39").Value = 0.75
$ws.Range("AA39").Value = 1.025
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = -1
$ws.Range("AD39").Value = 0.825

# Row 261
$ws.Range("B261").Value = 8157953
$ws.Range("E261").Value = "FK Jablonec"
$ws.Range("F261").Value = "Pardubice"
$ws.Range("G261").Value = 3
$ws.Range("H261").Value = 0
$ws.Range("K261").Value = "H"
$ws.Range("L261").Value = 1.833
$ws.Range("M261").Value = 3.75
$ws.Range("N261").Value = 4
$ws.Range("O261").Value = 1.727
$ws.Range("P261").Value = 3.8
$ws.Range("Q261").Value = 4.75
$ws.Range("R261").Value = -0.75
$ws.Range("S261").Value = 1.95
$ws.Range("T261").Value = 1.9
$ws.Range("X261").Value = 0.7270000000000001
$ws.Range("Y261").Value = -1
$ws.Range("AA261").Value = 0.95
$ws.Range("AB261").Value = -1
$ws.Range("AC261").Value = 0.875
$ws.Range("AD261").Value = -1

# Row 262
$ws.Range("B262").Value = 8157952
$ws.Range("E262").Value = "Bohemians 1905"
$ws.Range("F262").Value = "MFK Karvina"
$ws.Range("G262").Value = 1
$ws.Range("H262").Value = 3
$ws.Range("I262").Value = 1
$ws.Range("K262").Value = "A"
$ws.Range("M262").Value = 3.6
$ws.Range("N262").Value = 4.333
$ws.Range("O262").Value = 1.95
$ws.Range("P262").Value = 3.5
$ws.Range("Q262").Value = 3.9
$ws.Range("R262").Value = -0.5
$ws.Range("S262").Value = 1.9
$ws.Range("T262").Value = 1.95
$ws.Range("V262").Value = 1.85
$ws.Range("W262").Value = 2
$ws.Range("X262").Value = -1
$ws.Range("Z262").Value = 2.9
$ws.Range("AA262").Value = -1
$ws.Range("AB262").Value = 0.95
$ws.Range("AC262").Value = 0.8500000000000001

# Row 263
$ws.Range("B263").Value = 8157954
$ws.Range("E263").Value = "FC Trinity Zlin"
$ws.Range("F263").Value = "Ceske Budejovice"
$ws.Range("H263").Value = 1
$ws.Range("I263").Value = 0
$ws.Range("K263").Value = "D"
$ws.Range("L263").Value = 2.1
$ws.Range("M263").Value = 3.4
$ws.Range("N263").Value = 3.4
$ws.Range("O263").Value = 2.15
$ws.Range("P263").Value = 3.4
$ws.Range("Q263").Value = 3.3
$ws.Range("R263").Value = -0.25
$ws.Range("S263").Value = 1.85
$ws.Range("T263").Value = 2
$ws.Range("V263").Value = 1.875
$ws.Range("W263").Value = 1.975
$ws.Range("Y263").Value = 2.4
$ws.Range("Z263").Value = -1
$ws.Range("AA263").Value = -0.5
$ws.Range("AB263").Value = 0.5
$ws.Range("AC263").Value = -1
$ws.Range("AD263").Value = 0.9750000000000001
